# Update the Optical_Power sheet:
#  - Insert a new record row at row 26 (pushes existing rows 26-72 down to 27-73)
#  - Append two brand new record rows at the end (rows 74 and 75)
#
# Net effect matches the target diff: dimension grows from A1:P72 to A1:P75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at position 26, shifting current rows 26-72 down to 27-73.
$ws.Rows("26:26").Insert()

# Columns A, B, D, E and I hold values that look numeric (case numbers,
# dates written as text, comuna codes, OT numbers, flags) but must stay as
# plain text, matching the rest of the sheet. Force Text format first so
# Excel does not silently convert these strings into numbers/dates.
foreach ($col in @("A", "B", "D", "E", "I")) {
    $ws.Range($col + "26").NumberFormat = "@"
    $ws.Range($col + "74").NumberFormat = "@"
    $ws.Range($col + "75").NumberFormat = "@"
}

# 2) Populate the newly inserted row 26 with its data.
$row26 = @("5521", "4/8/2025", "EL PEREGRINO 3115", "11", "804569000", "Optical Power", "Pendiente", "Volvio a ingresar se inclino el poste - caso 6316", "1", "Aplomo", "Sin equipos", "Poste", -58.485232, -34.611573, "Devoto", "Capital Norte")
for ($i = 0; $i -lt $row26.Length; $i++) {
    $ws.Cells.Item(26, $i + 1).Value = $row26[$i]
}

# 3) Append two new rows (74 and 75) with brand new records.
$row74 = @("6330", "7/3/2025", "REPUBLICA DE LA INDIA 3106", "14", "807965776", "Optical Power", "Pendiente", "Picada e inclinada", "1", "Cambio", "Sin equipos", "Pasante", -58.413941, -34.57698, "Palermo", "Capital Sur")
for ($i = 0; $i -lt $row74.Length; $i++) {
    $ws.Cells.Item(74, $i + 1).Value = $row74[$i]
}

$row75 = @("-501", "7/3/2025", "Cabello 3107", "14", "807971967", "Optical Power", "Pendiente", "Aplomar", "0", "Aplomo", "Sin equipos", "Terminal", -58.405749, -34.58224, "Recoleta", "Capital Sur")
for ($i = 0; $i -lt $row75.Length; $i++) {
    $ws.Cells.Item(75, $i + 1).Value = $row75[$i]
}
